$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 237, shifting existing rows 237..347 down to 238..348.
$ws.Rows.Item(237).Insert()

# Populate the new row 237 with the new daily price record.
$ws.Range("A237").Value2 = 10
$ws.Range("B237").Value2 = "Vega Modelo de Temuco"
$ws.Range("C237").Value2 = "La Araucanía"
$ws.Range("D237").Value2 = 44460
$ws.Range("E237").Value2 = 9
$ws.Range("F237").Value2 = 100112045
$ws.Range("G237").Value2 = "Zapallo"
$ws.Range("H237").Value2 = "Paine"
$ws.Range("I237").Value2 = "1a (guarda)"
$ws.Range("J237").Value2 = 400
$ws.Range("K237").Value2 = 500
$ws.Range("L237").Value2 = 600
$ws.Range("M237").Value2 = 550
$ws.Range("N237").Value2 = "`$/kilo (volumen en unidades)"
$ws.Range("O237").Value2 = "Región del Maule"
$ws.Range("P237").Value2 = 550
$ws.Range("Q237").Value2 = 1
$ws.Range("R237").Value2 = "Hortaliza"

# Preserve the date-cell number format used throughout column D.
$ws.Range("D237").NumberFormat = $ws.Range("D238").NumberFormat
